$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 used to be "e편한세상두산1단지" (id 120178) and row 45 was
# "e편한세상두산2단지" (id 120179). The two complexes were consolidated into a
# single listing "e편한세상둔산" (id 182279, rich-text: "e" + formatted run).
# Update row 44 in place, then delete row 45 so every row below shifts up.

$ws.Cells.Item(44, 1).Value = 182279
$ws.Cells.Item(44, 2).Value = "e편한세상둔산"

$chars = $ws.Cells.Item(44, 2).Characters(2, 6)
$chars.Font.Name = "맑은 고딕"
$chars.Font.Size = 10

$ws.Rows.Item(45).Delete()
